$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0751464458737453
$ws.Range("D2").Value = 0.2119198634755605
$ws.Range("G2").Value = 0.1245324579833929
$ws.Range("H2").Value = 0.992
